# Natmi following Dr Hou advice
# Adds the "ECs" sending/target cluster to the Sema4d-Cd72 LR-pair sheet.
# "ECs" moves to the front of the cluster ordering (shared-string / category
# order becomes ECs, FAPs, M2, sCs) and the data grid is recomputed for the
# full 4x4 cluster cross-product (rows 2-17), growing the sheet from
# A1:T13 to A1:T17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-17 (columns A:T), in final order.
$rows = @(
  @("ECs","Sema4d","Cd72","ECs",2.0,0.6666666666666666,0.8730476666666668,2.619143,0.01740928848427011,0.01740928848427011,3.0,1.0,2.821352666666667,8.464058000000001,0.06564465474960138,0.06564465474960136,2.463175362477112,22.16857826229401,0.001142826731986122,0.001142826731986122),
  @("ECs","Sema4d","Cd72","FAPs",2.0,0.6666666666666666,0.8730476666666668,2.619143,0.01740928848427011,0.01740928848427011,3.0,1.0,4.044796666666667,12.13439,0.0941106313481093,0.09411063134810929,3.531300291974445,31.78170262777,0.001638399130576029,0.001638399130576029),
  @("ECs","Sema4d","Cd72","M2",2.0,0.6666666666666666,0.8730476666666668,2.619143,0.01740928848427011,0.01740928848427011,3.0,1.0,34.17538566666667,102.526157,0.7951616327615461,0.795161632761546,29.83674071371679,268.5306664234511,0.013843198256369,0.013843198256369),
  @("ECs","Sema4d","Cd72","sCs",2.0,0.6666666666666666,0.8730476666666668,2.619143,0.01740928848427011,0.01740928848427011,3.0,1.0,1.937633333333334,5.812900000000001,0.04508308114074334,0.04508308114074334,1.691646260522223,15.2248163447,0.0007848643653389579,0.0007848643653389579),
  @("FAPs","Sema4d","Cd72","ECs",3.0,1.0,1.252512666666667,3.757538,0.02497613266347325,0.02497613266347325,3.0,1.0,2.821352666666667,8.464058000000001,0.06564465474960138,0.06564465474960136,3.533779952133779,31.80401956920401,0.001639549605673944,0.001639549605673943),
  @("FAPs","Sema4d","Cd72","FAPs",3.0,1.0,1.252512666666667,3.757538,0.02497613266347325,0.02497613266347325,3.0,1.0,4.044796666666667,12.13439,0.0941106313481093,0.09411063134810929,5.066159059091111,45.59543153182,0.002350519613593603,0.002350519613593602),
  @("FAPs","Sema4d","Cd72","M2",3.0,1.0,1.252512666666667,3.757538,0.02497613266347325,0.02497613266347325,3.0,1.0,34.17538566666667,102.526157,0.7951616327615461,0.795161632761546,42.80510343571846,385.2459309214661,0.01986006242875638,0.01986006242875637),
  @("FAPs","Sema4d","Cd72","sCs",3.0,1.0,1.252512666666667,3.757538,0.02497613266347325,0.02497613266347325,3.0,1.0,1.937633333333334,5.812900000000001,0.04508308114074334,0.04508308114074334,2.426910293355556,21.8421926402,0.001126001015449335,0.001126001015449335),
  @("M2","Sema4d","Cd72","ECs",3.0,1.0,45.633473,136.900419,0.9099689814525027,0.9099689814525027,3.0,1.0,2.821352666666667,8.464058000000001,0.06564465474960138,0.06564465474960136,128.7481207378114,1158.733086640302,0.05973459962029596,0.05973459962029595),
  @("M2","Sema4d","Cd72","FAPs",3.0,1.0,45.633473,136.900419,0.9099689814525027,0.9099689814525027,3.0,1.0,4.044796666666667,12.13439,0.0941106313481093,0.09411063134810929,184.5781194788233,1661.20307530941,0.08563775535169099,0.08563775535169098),
  @("M2","Sema4d","Cd72","M2",3.0,1.0,45.633473,136.900419,0.9099689814525027,0.9099689814525027,3.0,1.0,34.17538566666667,102.526157,0.7951616327615461,0.795161632761546,1559.54153908442,14035.87385175978,0.7235724210541331,0.723572421054133),
  @("M2","Sema4d","Cd72","sCs",3.0,1.0,45.633473,136.900419,0.9099689814525027,0.9099689814525027,3.0,1.0,1.937633333333334,5.812900000000001,0.04508308114074334,0.04508308114074334,88.42093840056668,795.7884456051002,0.04102420542638276,0.04102420542638275),
  @("sCs","Sema4d","Cd72","ECs",3.0,1.0,2.389349666666667,7.168049000000001,0.04764559739975398,0.04764559739975399,3.0,1.0,2.821352666666667,8.464058000000001,0.06564465474960138,0.06564465474960136,6.741198053649113,60.67078248284201,0.003127678791645356,0.003127678791645355),
  @("sCs","Sema4d","Cd72","FAPs",3.0,1.0,2.389349666666667,7.168049000000001,0.04764559739975398,0.04764559739975399,3.0,1.0,4.044796666666667,12.13439,0.0941106313481093,0.09411063134810929,9.664433567234445,86.97990210511001,0.004483957252248682,0.004483957252248682),
  @("sCs","Sema4d","Cd72","M2",3.0,1.0,2.389349666666667,7.168049000000001,0.04764559739975398,0.04764559739975399,3.0,1.0,34.17538566666667,102.526157,0.7951616327615461,0.795161632761546,81.6569463508548,734.9125171576932,0.03788595102228765,0.03788595102228765),
  @("sCs","Sema4d","Cd72","sCs",3.0,1.0,2.389349666666667,7.168049000000001,0.04764559739975398,0.04764559739975399,3.0,1.0,1.937633333333334,5.812900000000001,0.04508308114074334,0.04508308114074334,4.629683559122223,41.66715203210001,0.002148010333572299,0.002148010333572299)
)

$nrows = $rows.Count
$ncols = $rows[0].Count
$data = New-Object 'object[,]' $nrows, $ncols
for ($i = 0; $i -lt $nrows; $i++) {
    for ($j = 0; $j -lt $ncols; $j++) {
        $data[$i, $j] = $rows[$i][$j]
    }
}

$ws.Range("A2:T17").Value = $data
